$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update facility limit (E) and facility utilisation (F) values for rows 2-6
$ws.Range("E2").Value = 800000
$ws.Range("F2").Value = 761579.37

$ws.Range("E3").Value = 800000
$ws.Range("F3").Value = 761579.37

$ws.Range("E4").Value = 800000
$ws.Range("F4").Value = 761579.37

$ws.Range("E5").Value = 800000
$ws.Range("F5").Value = 761579.37
$ws.Range("G5").Value = 456
$ws.Range("H5").Value = 3938753.8

$ws.Range("E6").Value = 800000
$ws.Range("G6").Value = 761579.37

# Make column G match the width/bestFit of columns E:F
$ws.Range("G1").EntireColumn.ColumnWidth = 16.33203125

# Update the active selection shown in the sheet view
$ws.Range("E2:H6").Select()
